$d = $word.ActiveDocument

# The "Ingredientes" heading paragraph (2nd paragraph in the document).
$ingredientesPara = $d.Paragraphs(2)
$ingredientesRange = $ingredientesPara.Range

# Split off a throw-away paragraph right after it. Word clones the
# paragraph mark's formatting (pPr/rPr) onto the new, empty paragraph, so
# anything typed into it immediately picks up the same run formatting
# (Century Gothic / F5886B / 30pt) as "Ingredientes" without us having to
# respecify every rFonts sub-font by hand.
$ingredientesRange.InsertParagraphAfter()
$scratchPara = $d.Paragraphs(3)
$scratchRange = $scratchPara.Range
$scratchRange.InsertAfter(":")

# Grab just the freshly typed ":" (with its inherited formatting) and copy
# it onto the clipboard.
$colonStart = $scratchPara.Range.Start
$colonRange = $d.Range($colonStart, $colonStart + 1)
$colonRange.Copy()

# Paste it in as a new run immediately after "Ingredientes", i.e. right
# before that paragraph's own paragraph mark. Pasting like this appends a
# sibling <w:r> instead of merging text into the existing "Ingredientes"
# run, and - crucially - it keeps the "Ingredientes" paragraph's own
# identity (paraId / rsid attributes) untouched, unlike deleting a
# paragraph mark to splice paragraphs back together would.
$pasteAt = $ingredientesPara.Range.End - 1
$pasteRange = $d.Range($pasteAt, $pasteAt)
$pasteRange.Paste()

# Clean up the scratch paragraph we used as a formatting donor - its
# paragraph index shifted by one once the paragraph mark split, so
# re-fetch it by index rather than reusing the old object.
$scratchParaAgain = $d.Paragraphs(3)
$scratchDeleteRange = $d.Range($scratchParaAgain.Range.Start, $scratchParaAgain.Range.End)
$scratchDeleteRange.Delete()
